$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (row 1 is the header row), shifting all
# existing data rows down by one.
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits formatting from the header row above it
# (bold/centered). Clear that so it matches the plain formatting used by the
# other data rows, then re-apply the date/time number format used in column C.
$ws.Range("A2:D2").ClearFormats()
$ws.Range("C2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row's values (column B / "comment" stays blank).
$ws.Range("A2").Value = 5
$ws.Range("C2").Value = 46034.61868471065
$ws.Range("D2").Value = "OWY1NGRiM2EtMjUxMy00YzZjLTg0ZDctMTBhZGU5MjQzZjY2OjU3MDE2"
